# Generate Report for Handback
# Updates the localization status workbook to reflect that the
# "1c5008db-8f21-4ef3-ae9e-8bfa997a412b" file has now been handed back
# (was previously "Ready for handoff"), and stamps the new handback
# timestamps for both locales.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G2").Value = "2016-03-03 13:23:29"
$wsZhCn.Range("G3").Value = "2016-03-03 13:23:29"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G2").Value = "2016-03-03 13:23:54"
$wsDeDe.Range("G3").Value = "2016-03-03 13:23:54"
